$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Burnup data table (dates shift by one working day, Total
#     drops from 20 to 14, and the "Completed" series is populated) ---
$ws.Range("A2").Value = 43031
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = 43034
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = 43037
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = 43040
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = 43042
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 6

$ws.Range("A7").Value = 43044
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 8

# --- Chart formatting: un-bold the chart title and the axis tick labels ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart

$chart.ChartTitle.Font.Bold = $false

$catAxis = $chart.Axes(1)
$catAxis.TickLabels.Font.Bold = $false

$valAxis = $chart.Axes(2)
$valAxis.TickLabels.Font.Bold = $false
